# Apply corrected timetable assignments to Section_A and Section_B sheets
# (fix Excel generation for timetables: reshuffled course slots).
$wb = $excel.ActiveWorkbook

$wsA = $wb.Worksheets.Item("Section_A")
$wsA.Range("B2").Value = "CS309"
$wsA.Range("C2").Value = "CS309"
$wsA.Range("D2").Value = "Free"
$wsA.Range("E2").Value = "CS309"
$wsA.Range("B3").Value = "Free"
$wsA.Range("C3").Value = "CS303"
$wsA.Range("F3").Value = "CS303"
$wsA.Range("B5").Value = "Free"
$wsA.Range("C5").Value = "CS461"
$wsA.Range("D5").Value = "CS461"
$wsA.Range("E5").Value = "CS303"
$wsA.Range("C6").Value = "CS304"
$wsA.Range("E6").Value = "Free"
$wsA.Range("B7").Value = "CS461"
$wsA.Range("D7").Value = "CS304"
$wsA.Range("F7").Value = "CS304"

$wsB = $wb.Worksheets.Item("Section_B")
$wsB.Range("C2").Value = "CS309"
$wsB.Range("D2").Value = "Free"
$wsB.Range("E2").Value = "CS304"
$wsB.Range("F2").Value = "Free"
$wsB.Range("B3").Value = "CS309"
$wsB.Range("C3").Value = "CS304"
$wsB.Range("D3").Value = "Free"
$wsB.Range("E3").Value = "CS309"
$wsB.Range("F3").Value = "CS461"
$wsB.Range("C5").Value = "CS461"
$wsB.Range("D5").Value = "CS461"
$wsB.Range("E5").Value = "Free"
$wsB.Range("F5").Value = "Free"
$wsB.Range("B6").Value = "CS461"
$wsB.Range("C6").Value = "CS303"
$wsB.Range("D6").Value = "Free"
$wsB.Range("F6").Value = "CS303"
$wsB.Range("B7").Value = "CS304"
$wsB.Range("C7").Value = "Free"
$wsB.Range("D7").Value = "CS303"
$wsB.Range("F7").Value = "Free"

Write-Host "Timetable updated"
